$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Match the formatting already used by the other data rows (e.g. row 22)
# before filling in the new values, by copying just the cell formatting.
$ws.Range("A22:H22").Copy()
$ws.Range("A23:H23").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Add the 2023 APWH exam score row.
$ws.Range("A23").Value = 2023
$ws.Range("B23").Value = 15
$ws.Range("C23").Value = 22
$ws.Range("D23").Value = 28
$ws.Range("E23").Value = 22
$ws.Range("F23").Value = 13
$ws.Range("G23").Value = 3.04
$ws.Range("H23").Value = 356000
